$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 into I1:J1 (reuses the same style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Set new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill columns I (I0, always 1) and J (IF, mirrors column H / IP) for rows 2-28
for ($r = 2; $r -le 28; $r++) {
    $ip = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ip
}
